$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.584.52'
$ws.Range("E2").Value = '  +4.11%  '
$ws.Range("D3").Value = '1.743.92'
$ws.Range("E3").Value = '  +4.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.94%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4821'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.20%  '
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06263'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.73%  '
$ws.Range("D10").Value = '1.742.85'
$ws.Range("E10").Value = '  +4.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07143'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6225'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.522'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '26.582.90'
$ws.Range("E17").Value = '  +4.12%  '
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006908'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.48%  '
$ws.Range("D21").Value = '1.968.14'
$ws.Range("E21").Value = '  +4.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.616'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.874'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.378'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.815'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.426'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.745'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.27%  '
$ws.Range("E32").Value = '  +1.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04585'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.616'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6389'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9391'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '113.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.69%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.985'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.27%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.422'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.81%  '
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.775'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +17.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01516'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3923'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.788'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.04%  '
$ws.Range("E46").Value = '  +9.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05334'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.945'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.267'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3451'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.27%  '
